# Update to GAS hierarchy
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "Outcome" values down one level in the GAS hierarchy
$ws.Range("B2").Value = 40
$ws.Range("B3").Value = 30
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 10
$ws.Range("B6").Value = 0

# Widen column A slightly to fit the labels
$ws.Columns("A").ColumnWidth = 23.5

# Update the current selection
$ws.Range("C2:C6").Select()
